$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(
    '24 19:11>>> CFD893A460   Gary Tsai',
    '24 19:11>>> 0FD88F2580   Carolyn Alana',
    '30 16:43>>> 4FD8A85BA0   hojin euam',
    '30 16:43>>> 4FD8A85BA0   hojin euam',
    '30 16:43>>> CFD8AF34C0   Megan Eng',
    '30 16:47>>> 0FD8B42CC0   ',
    '30 16:52>>> 0FD8A290A0   Paul Fabro',
    '30 17:17>>> 8FD89F60C0    ',
    '30 17:32>>> 0FD8AE0B20   Anthony Rivera',
    '30 17:39>>> 0FD8AF13C0   Rondell Holland',
    '30 17:39>>> 0FD8A6BA00   John Orellana',
    '30 17:39>>> 0FD8AFD720   Jennifer Mercedes',
    '30 17:39>>> CFD8AFA4C0    ',
    '30 17:40>>> 8FD8A3A820   ',
    '30 17:40>>> CFD8AFA4C0    ',
    '30 17:40>>> CFD8AFA4C0   arango juan',
    '30 17:40>>> CFD8A6B7E0   peter sormilic',
    '30 17:41>>> 8FD8A3A820    ',
    '30 17:41>>> 8FD8A3A820   cowell,truman',
    '30 17:41>>> 0FD8B42CC0    ',
    '30 17:41>>> 0FD8B42CC0   amin shahid',
    '30 17:52>>> 8FD8AC5D40   Rudelina Calcano',
    '30 17:53>>> CFD8A948E0   Vanessa Cordero',
    '30 17:53>>> CFD8A948E0   Vanessa Cordero',
    '30 17:53>>> CFD8A948E0   Vanessa Cordero',
    '30 17:53>>> 8FD8AC5D40   Rudelina Calcano',
    '30 17:57>>> 8FD8B367A0   ',
    '30 18:01>>> 4FD8B36A40   ',
    '30 18:01>>> 4FD8B36A40    ',
    '30 18:31>>> 0FD88F2580   Carolyn Alana',
    '30 19:20>>> 0FD8AFD720   Jennifer Mercedes',
    '30 19:21>>> 0FD8AFD720   Jennifer Mercedes',
    '30 19:54>>> 0FD8B45C20   ',
    '31 08:27>>> 8FD8B68DE0   ',
    '31 08:27>>> 0FD8A6BA00   John Orellana',
    '31 08:27>>> 8FD8B68DE0    ',
    '31 08:28>>> 8FD8B68DE0   Miguel Martillo',
    '31 08:52>>> CFD88517C0   Joseph Giordano',
    '31 09:01>>> 0FD8A61420   Duban hershely',
    '31 09:01>>> CFD8A3DFA0   eugene marmontov',
    '31 09:20>>> 0FD8AFCD00   Joshua Itteera',
    '31 09:21>>> 0FD8AFCD00   Joshua Itteera',
    '31 09:26>>> 0FD8B60440    ',
    '31 09:54>>> 0FD89CD000   David Schachner',
    '31 10:00>>> CFD8ADF120   Patrick Mclaughlin',
    '31 10:20>>> 0FD8AD9C60   everton mendozaholmes',
    '31 10:34>>> 8FD87CA8C0   Dwayne Stallworth',
    '31 10:41>>> 4FD8AE2B20   Brandon haden',
    '31 10:51>>> 8FD8A60A60   ferguson elaina',
    '31 10:53>>> 0FD8AE0B20   Anthony Rivera',
    '31 11:45>>> 4FD8A51080   Calvin Y Au',
    '31 12:02>>> 4FD8B5DC60   Nadine Espinosa',
    '31 13:15>>> 8FD8B65040   ',
    '31 13:27>>> 0FD8B60440    ',
    '31 13:37>>> CFD8AEC840   ',
    '31 13:38>>> 4FD8A33DE0   Anne Crosby',
    '31 13:43>>> 0FD8A9BD80   Jaspreet Kaur',
    '31 13:45>>> CFD8A948E0   Vanessa Cordero',
    '31 13:45>>> CFD8A948E0   Vanessa Cordero',
    '31 13:48>>> 8FD8ADBD20   ',
    '31 13:56>>> 0FD8AF13C0   Rondell Holland',
    '31 14:07>>> CFD893A460   Gary Tsai',
    '31 14:08>>> 8FD8A91340   Jesse Silkworth',
    '31 14:18>>> CFD88517C0   Joseph Giordano',
    '31 14:20>>> 0FD8AFD720   Jennifer Mercedes'
)

$startRow = 177
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $values[$i]
}
